$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four data rows (2-5) get their values permuted (dates reordered along
# with their associated Calidad/Volumen/Precio/Unidad/Kg-Unidades columns).
# New row 2 <- old row 4
# New row 3 <- old row 2
# New row 4 <- old row 5
# New row 5 <- old row 3

$newRow2 = @(44623, "Primera", 300, 1800, 2000, 1900, "`$/paquete", 1900, 1)
$newRow3 = @(44370, "Segunda", 100, 1000, 1200, 1080, "`$/docena de matas", 180, 6)
$newRow4 = @(44267, "Primera", 120, 1500, 1800, 1650, "`$/docena de matas", 275, 6)
$newRow5 = @(44377, "Segunda", 550, 2000, 2800, 2364, "`$/docena de matas", 394, 6)

$rows = @{ 2 = $newRow2; 3 = $newRow3; 4 = $newRow4; 5 = $newRow5 }

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("I$r").Value = $vals[1]
    $ws.Range("J$r").Value = $vals[2]
    $ws.Range("K$r").Value = $vals[3]
    $ws.Range("L$r").Value = $vals[4]
    $ws.Range("M$r").Value = $vals[5]
    $ws.Range("N$r").Value = $vals[6]
    $ws.Range("P$r").Value = $vals[7]
    $ws.Range("Q$r").Value = $vals[8]
}
